$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7248245.5
$ws.Range("I40").Value = 1873.8334
$ws.Range("J40").Value = 55557388
$ws.Range("K40").Value = 1873.8334
$ws.Range("L40").Value = 55557388
$ws.Range("M40").Value = -1698.8334
$ws.Range("N40").Value = -55557738

$ws.Range("H86").Value = 4390.6665
$ws.Range("I86").Value = 4554.4375
$ws.Range("K86").Value = 4554.4375
$ws.Range("M86").Value = -3431.4375

$ws.Range("H89").Value = 4390.6665
$ws.Range("I89").Value = 4554.4375
$ws.Range("K89").Value = 22772.1875
$ws.Range("M89").Value = -17156.1875

$ws.Range("H98").Value = 998.5
$ws.Range("I98").Value = 641.5333000000001
$ws.Range("J98").Value = 1763.4286
$ws.Range("K98").Value = 641.5333000000001
$ws.Range("L98").Value = 1763.4286
$ws.Range("M98").Value = 856.4666999999999
$ws.Range("N98").Value = -4759.4286

$ws.Range("H112").Value = 21979522
$ws.Range("I112").Value = 466.66666
$ws.Range("J112").Value = 24846356
$ws.Range("K112").Value = 1399.99998
$ws.Range("L112").Value = 74539068
$ws.Range("M112").Value = -291.9999800000001
$ws.Range("N112").Value = -74541284

$ws.Range("H122").Value = 998.5
$ws.Range("I122").Value = 641.5333000000001
$ws.Range("J122").Value = 1763.4286
$ws.Range("K122").Value = 1924.5999
$ws.Range("L122").Value = 5290.2858
$ws.Range("M122").Value = 525.4000999999998
$ws.Range("N122").Value = -10190.2858

$ws.Range("H137").Value = 2081.8708
$ws.Range("I137").Value = 1578.2693
$ws.Range("J137").Value = 4700.6
$ws.Range("K137").Value = 4734.8079
$ws.Range("L137").Value = 14101.8
$ws.Range("M137").Value = -2184.8079
$ws.Range("N137").Value = -19201.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5795.33
$ws.Range("I32").Value = 4393.809
$ws.Range("J32").Value = 17134.908
$ws.Range("K32").Value = 4393.809
$ws.Range("L32").Value = 17134.908
$ws.Range("M32").Value = -4106.809
$ws.Range("N32").Value = -17708.908

$ws.Range("H45").Value = 175893.33
$ws.Range("I45").Value = 263177.5
$ws.Range("J45").Value = 1325
$ws.Range("K45").Value = 263177.5
$ws.Range("L45").Value = 1325
$ws.Range("M45").Value = -262800.5
$ws.Range("N45").Value = -2079

$ws.Range("H61").Value = 208393.89
$ws.Range("I61").Value = 4992.161
$ws.Range("J61").Value = 558696.9
$ws.Range("K61").Value = 4992.161
$ws.Range("L61").Value = 558696.9
$ws.Range("M61").Value = -4780.161
$ws.Range("N61").Value = -559120.9

$ws.Range("H88").Value = 600
$ws.Range("J88").Value = 600
$ws.Range("L88").Value = 600
$ws.Range("N88").Value = -1412

$ws.Range("H91").Value = 600
$ws.Range("J91").Value = 600
$ws.Range("L91").Value = 600
$ws.Range("N91").Value = -3408

$ws.Range("H110").Value = 13567.52
$ws.Range("I110").Value = 17352
$ws.Range("J110").Value = 1583.3334
$ws.Range("K110").Value = 17352
$ws.Range("L110").Value = 1583.3334
$ws.Range("M110").Value = -15307
$ws.Range("N110").Value = -5673.3334

$ws.Range("H136").Value = 208393.89
$ws.Range("I136").Value = 4992.161
$ws.Range("J136").Value = 558696.9
$ws.Range("K136").Value = 14976.483
$ws.Range("L136").Value = 1676090.7
$ws.Range("M136").Value = -12426.483
$ws.Range("N136").Value = -1681190.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 202250.06
$ws.Range("J107").Value = 1535.75
$ws.Range("L107").Value = 1535.75
$ws.Range("N107").Value = -5375.75

$ws.Range("H134").Value = 18143.516
$ws.Range("I134").Value = 3829.453
$ws.Range("J134").Value = 68719.87
$ws.Range("K134").Value = 11488.359
$ws.Range("L134").Value = 206159.61
$ws.Range("M134").Value = -8953.359
$ws.Range("N134").Value = -211229.61

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 218.4
$ws.Range("I22").Value = 197
$ws.Range("J22").Value = 250.5
$ws.Range("K22").Value = 197
$ws.Range("L22").Value = 250.5
$ws.Range("M22").Value = 153
$ws.Range("N22").Value = -950.5

$ws.Range("H132").Value = 8698772
$ws.Range("I132").Value = 20002204
$ws.Range("J132").Value = 3824.3076
$ws.Range("K132").Value = 60006612
$ws.Range("L132").Value = 11472.9228
$ws.Range("M132").Value = -60004082
$ws.Range("N132").Value = -16532.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H131").Value = 2858574.2
$ws.Range("I131").Value = 8334041.5
$ws.Range("J131").Value = 1808.5217
$ws.Range("K131").Value = 25002124.5
$ws.Range("L131").Value = 5425.5651
$ws.Range("M131").Value = -24997084.5
$ws.Range("N131").Value = -15505.5651

$ws.Range("H134").Value = 15992.333
$ws.Range("I134").Value = 18988.5
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 56965.5
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -51895.5
$ws.Range("N134").Value = -40140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 97489330
$ws.Range("I122").Value = 177470770
$ws.Range("J122").Value = 37503252
$ws.Range("K122").Value = 532412310
$ws.Range("L122").Value = 112509756
$ws.Range("M122").Value = -532409860
$ws.Range("N122").Value = -112514656

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 31637.715
$ws.Range("J24").Value = 31637.715
$ws.Range("L24").Value = 31637.715
$ws.Range("N24").Value = -32323.715

$ws.Range("H40").Value = 90911816
$ws.Range("I40").Value = 125002560
$ws.Range("J40").Value = 3151.6667
$ws.Range("K40").Value = 125002560
$ws.Range("L40").Value = 3151.6667
$ws.Range("M40").Value = -125002424
$ws.Range("N40").Value = -3423.6667

$ws.Range("H55").Value = 103.8
$ws.Range("I55").Value = 131
$ws.Range("J55").Value = 63
$ws.Range("K55").Value = 131
$ws.Range("L55").Value = 63
$ws.Range("M55").Value = 42
$ws.Range("N55").Value = -409

$ws.Range("H122").Value = 5434628
$ws.Range("I122").Value = 5959285
$ws.Range("J122").Value = 3336000
$ws.Range("K122").Value = 17877855
$ws.Range("L122").Value = 10008000
$ws.Range("M122").Value = -17875405
$ws.Range("N122").Value = -10012900

$ws.Range("H136").Value = 9866.484
$ws.Range("I136").Value = 7372.5454
$ws.Range("J136").Value = 14854.363
$ws.Range("K136").Value = 22117.6362
$ws.Range("L136").Value = 44563.089
$ws.Range("M136").Value = -19567.6362
$ws.Range("N136").Value = -49663.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1331.5927
$ws.Range("I132").Value = 750.29266
$ws.Range("J132").Value = 3164.923
$ws.Range("K132").Value = 2250.87798
$ws.Range("L132").Value = 9494.769
$ws.Range("M132").Value = 279.1220200000002
$ws.Range("N132").Value = -14554.769

$ws.Range("H136").Value = 3089526.8
$ws.Range("I136").Value = 3396.7083
$ws.Range("J136").Value = 5558430.5
$ws.Range("K136").Value = 10190.1249
$ws.Range("L136").Value = 16675291.5
$ws.Range("M136").Value = -7640.124899999999
$ws.Range("N136").Value = -16680391.5
